$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'63.488.87"
$ws.Range("E2").Value = "  +0.97%  "
$ws.Range("D3").Value = "'3.098.94"
$ws.Range("E3").Value = "  -0.44%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'583.35"
$ws.Range("E5").Value = "  -0.06%  "
$ws.Range("D6").Value = "'144.77"
$ws.Range("E6").Value = "  -0.17%  "
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("D8").Value = "'3.090.83"
$ws.Range("E8").Value = "  -0.48%  "
$ws.Range("D9").Value = "'0.528"
$ws.Range("E9").Value = "  -0.34%  "
$ws.Range("E10").Value = "  +6.66%  "
$ws.Range("D11").Value = "'5.61"
$ws.Range("E11").Value = "  -2.58%  "
$ws.Range("E12").Value = "  -2.43%  "
$ws.Range("D13").Value = "'0.0000246"
$ws.Range("E13").Value = "  -1.14%  "
$ws.Range("D14").Value = "'37.24"
$ws.Range("E14").Value = "  +4.60%  "
$ws.Range("E15").Value = "  -1.17%  "
$ws.Range("D16").Value = "'3.611.93"
$ws.Range("E16").Value = "  -0.42%  "
$ws.Range("D17").Value = "'63.353.12"
$ws.Range("E17").Value = "  +0.90%  "
$ws.Range("E18").Value = "  -1.03%  "
$ws.Range("D19").Value = "'3.094.40"
$ws.Range("E19").Value = "  -0.59%  "
$ws.Range("D20").Value = "'461.83"
$ws.Range("E20").Value = "  -1.41%  "
$ws.Range("D21").Value = "'14.23"
$ws.Range("E21").Value = "  +0.99%  "
$ws.Range("E22").Value = "  -0.52%  "
$ws.Range("E23").Value = "  -1.20%  "
$ws.Range("D24").Value = "'81.27"
$ws.Range("E24").Value = "  -0.91%  "
$ws.Range("D25").Value = "'12.91"
$ws.Range("E25").Value = "  -2.96%  "
$ws.Range("E26").Value = "  -1.56%  "
$ws.Range("E27").Value = "  +0.10%  "
$ws.Range("D28").Value = "'8.98"
$ws.Range("E28").Value = "  +8.50%  "
$ws.Range("E29").Value = "  +0.04%  "
$ws.Range("E30").Value = "  -0.34%  "
$ws.Range("E31").Value = "  -2.14%  "
$ws.Range("D32").Value = "'6.85"
$ws.Range("E32").Value = "  +0.53%  "
$ws.Range("E33").Value = "  -0.81%  "
$ws.Range("E34").Value = "  -1.34%  "
$ws.Range("E35").Value = "  -2.48%  "
$ws.Range("D36").Value = "'3.40"
$ws.Range("E36").Value = "  +3.83%  "
$ws.Range("E37").Value = "  -0.79%  "
$ws.Range("E38").Value = "  -3.47%  "
$ws.Range("E39").Value = "  -0.98%  "
$ws.Range("D40").Value = "'50.34"
$ws.Range("E40").Value = "  -1.14%  "
$ws.Range("D41").Value = "'434.26"
$ws.Range("E41").Value = "  +0.49%  "
$ws.Range("D42").Value = "'8.71"
$ws.Range("E42").Value = "  -0.34%  "
$ws.Range("E43").Value = "  -0.24%  "
$ws.Range("D44").Value = "'2.879.55"
$ws.Range("E44").Value = "  -1.75%  "
$ws.Range("E45").Value = "  -3.25%  "
$ws.Range("E46").Value = "  -2.92%  "
$ws.Range("D47").Value = "'35.85"
$ws.Range("E47").Value = "  +0.94%  "
$ws.Range("D49").Value = "'123.30"
$ws.Range("E50").Value = "  -1.40%  "
$ws.Range("D51").Value = "'24.07"
$ws.Range("E51").Value = "  -2.23%  "

Write-Output "done"
